$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - set text for all header cells (A1:O1)
$ws.Range("A1").Value = "#_Agents"
$ws.Range("B1").Value = "Coverage"
$ws.Range("C1").Value = "Avg_Total_Rounds"
$ws.Range("D1").Value = "Avg_Expl_Cost"
$ws.Range("E1").Value = "Avg_Expl_Eff"
$ws.Range("F1").Value = "Avg_Round_Time"
$ws.Range("G1").Value = "Avg_Agent_Step_Time"
$ws.Range("H1").Value = "Avg_Experiment_Time"
$ws.Range("I1").Value = "Std_Total_Rounds"
$ws.Range("J1").Value = "Std_Expl_Cost"
$ws.Range("K1").Value = "Std_Expl_Eff"
$ws.Range("L1").Value = "Std_Round_Time"
$ws.Range("M1").Value = "Std_Agent_Step_Time"
$ws.Range("N1").Value = "Std_Experiment_Time"
$ws.Range("O1").Value = "Obs_Prob"

# Apply the existing header style (s="1", same as A1:H1) to the newly
# added header cells I1:O1 so they match the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 56.326
$ws.Range("D2").Value = 56.326
$ws.Range("E2").Value = 3.06145044
$ws.Range("F2").Value = 0.17942226
$ws.Range("G2").Value = 0.17942226
$ws.Range("H2").Value = 9.873369720000001
$ws.Range("I2").Value = 7.890770083142077
$ws.Range("J2").Value = 7.890770083142077
$ws.Range("K2").Value = 0.4448968254753483
$ws.Range("L2").Value = 0.03238795169074209
$ws.Range("M2").Value = 0.03238795169074209
$ws.Range("N2").Value = 0.7670348691208346
$ws.Range("O2").Value = 0.15

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 90.764
$ws.Range("D3").Value = 90.764
$ws.Range("E3").Value = 1.89724074
$ws.Range("F3").Value = 0.12831342
$ws.Range("G3").Value = 0.12831342
$ws.Range("H3").Value = 11.45746828
$ws.Range("I3").Value = 12.51277119129777
$ws.Range("J3").Value = 12.51277119129777
$ws.Range("K3").Value = 0.2613756292341876
$ws.Range("L3").Value = 0.01917586704839485
$ws.Range("M3").Value = 0.01917586704839485
$ws.Range("N3").Value = 1.028636595057762
$ws.Range("O3").Value = 0.85

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 28.206
$ws.Range("D4").Value = 56.384
$ws.Range("E4").Value = 3.11549184
$ws.Range("F4").Value = 0.26137636
$ws.Range("G4").Value = 0.13068832
$ws.Range("H4").Value = 3.58765286
$ws.Range("I4").Value = 5.528534470071933
$ws.Range("J4").Value = 11.05209670982691
$ws.Range("K4").Value = 0.6254129187627363
$ws.Range("L4").Value = 0.04850405475633983
$ws.Range("M4").Value = 0.02425202117432028
$ws.Range("N4").Value = 0.4988782814203945
$ws.Range("O4").Value = 0.15

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 48.284
$ws.Range("D5").Value = 94.596
$ws.Range("E5").Value = 1.83554956
$ws.Range("F5").Value = 0.18118334
$ws.Range("G5").Value = 0.0905917
$ws.Range("H5").Value = 4.26867266
$ws.Range("I5").Value = 8.54760006866651
$ws.Range("J5").Value = 15.6988436320248
$ws.Range("K5").Value = 0.3042159677766882
$ws.Range("L5").Value = 0.03365826674808688
$ws.Range("M5").Value = 0.01682925519953698
$ws.Range("N5").Value = 0.5337443728023
$ws.Range("O5").Value = 0.85

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 14.708
$ws.Range("D6").Value = 58.784
$ws.Range("E6").Value = 3.05158964
$ws.Range("F6").Value = 0.302603
$ws.Range("G6").Value = 0.07565064
$ws.Range("H6").Value = 1.0814634
$ws.Range("I6").Value = 3.574954068682054
$ws.Range("J6").Value = 14.30967776560713
$ws.Range("K6").Value = 0.7588076778170465
$ws.Range("L6").Value = 0.07525947316852016
$ws.Range("M6").Value = 0.01881495149049232
$ws.Range("N6").Value = 0.2729347769959476
$ws.Range("O6").Value = 0.15

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 25.214
$ws.Range("D7").Value = 93.258
$ws.Range("E7").Value = 1.86439348
$ws.Range("F7").Value = 0.234371
$ws.Range("G7").Value = 0.05859298
$ws.Range("H7").Value = 1.44500108
$ws.Range("I7").Value = 5.087235391433712
$ws.Range("J7").Value = 16.07508172623113
$ws.Range("K7").Value = 0.3123801056987501
$ws.Range("L7").Value = 0.04714704591463768
$ws.Range("M7").Value = 0.01178696204575786
$ws.Range("N7").Value = 0.2722654843412229
$ws.Range("O7").Value = 0.85

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 9.436
$ws.Range("D8").Value = 56.512
$ws.Range("E8").Value = 3.19547968
$ws.Range("F8").Value = 0.34213698
$ws.Range("G8").Value = 0.05702296
$ws.Range("H8").Value = 0.5271577999999999
$ws.Range("I8").Value = 2.590622831155978
$ws.Range("J8").Value = 15.52365368227772
$ws.Range("K8").Value = 0.8016572794298218
$ws.Range("L8").Value = 0.09577927735277626
$ws.Range("M8").Value = 0.01596313812494248
$ws.Range("N8").Value = 0.1780111688192873
$ws.Range("O8").Value = 0.15

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 0.99992
$ws.Range("C9").Value = 17.204
$ws.Range("D9").Value = 88.392
$ws.Range("E9").Value = 1.96636592
$ws.Range("F9").Value = 0.25559048
$ws.Range("G9").Value = 0.04259843999999999
$ws.Range("H9").Value = 0.720108
$ws.Range("I9").Value = 4.034717076049761
$ws.Range("J9").Value = 14.69483625035848
$ws.Range("K9").Value = 0.3371963218867333
$ws.Range("L9").Value = 0.05489555714407055
$ws.Range("M9").Value = 0.009149152394921126
$ws.Range("N9").Value = 0.1813006339595874
$ws.Range("O9").Value = 0.85

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 6.984
$ws.Range("D10").Value = 55.722
$ws.Range("E10").Value = 3.24341424
$ws.Range("F10").Value = 0.32172954
$ws.Range("G10").Value = 0.04021635999999999
$ws.Range("H10").Value = 0.28002288
$ws.Range("I10").Value = 1.946103449787433
$ws.Range("J10").Value = 15.59587055462802
$ws.Range("K10").Value = 0.8187663575267718
$ws.Range("L10").Value = 0.08949062761681192
$ws.Range("M10").Value = 0.01118623490066724
$ws.Range("N10").Value = 0.1115522876622032
$ws.Range("O10").Value = 0.15

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 13.53
$ws.Range("D11").Value = 83.808
$ws.Range("E11").Value = 2.09872782
$ws.Range("F11").Value = 0.24704478
$ws.Range("G11").Value = 0.03088046
$ws.Range("H11").Value = 0.41095022
$ws.Range("I11").Value = 3.901736420388796
$ws.Range("J11").Value = 16.74885621076071
$ws.Range("K11").Value = 0.4305364061613236
$ws.Range("L11").Value = 0.05288329339585787
$ws.Range("M11").Value = 0.006610008442628742
$ws.Range("N11").Value = 0.1284208122954096
$ws.Range("O11").Value = 0.85

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 5.598
$ws.Range("D12").Value = 55.802
$ws.Range("E12").Value = 3.22291532
$ws.Range("F12").Value = 0.29306114
$ws.Range("G12").Value = 0.02930602
$ws.Range("H12").Value = 0.16705926
$ws.Range("I12").Value = 1.425644037564668
$ws.Range("J12").Value = 14.20896143883643
$ws.Range("K12").Value = 0.8125884086726292
$ws.Range("L12").Value = 0.08753907072334395
$ws.Range("M12").Value = 0.008753819730800876
$ws.Range("N12").Value = 0.07737972920757008
$ws.Range("O12").Value = 0.15

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 10.83
$ws.Range("D13").Value = 75.244
$ws.Range("E13").Value = 2.3328834
$ws.Range("F13").Value = 0.23279148
$ws.Range("G13").Value = 0.02327926
$ws.Range("H13").Value = 0.25024616
$ws.Range("I13").Value = 3.355535611266397
$ws.Range("J13").Value = 15.03117881848748
$ws.Range("K13").Value = 0.451108978007441
$ws.Range("L13").Value = 0.04479983438508968
$ws.Range("M13").Value = 0.004480113490322801
$ws.Range("N13").Value = 0.08712085956875308
$ws.Range("O13").Value = 0.85

